$wb = $excel.ActiveWorkbook

# --- Existing sheets (by name, to be safe regardless of index order) ---
$wsLogin  = $wb.Worksheets.Item("LoginPage")
$wsReset  = $wb.Worksheets.Item("ResetPage")
$wsHome   = $wb.Worksheets.Item("HomePage")
$wsUsers  = $wb.Worksheets.Item("UsersPage")

# --- Populate new data on the existing "UsersPage" sheet -----------------
# Cell write order controls shared-string allocation order, so B-then-A
# within each new row to match the target string table layout.
$wsUsers.Range("B2").Value = "User added successfully"
$wsUsers.Range("A2").Value = "AddUserSuccessMessage"
$wsUsers.Range("B3").Value = "No matching records found"
$wsUsers.Range("A3").Value = "SearchWithInvalidDataMessage"

# --- Add the two new worksheets at the end of the workbook ---------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsAddUser = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsAddUser.Name = "AddUsersPage"
$wsAddUser.Range("A1").Value = "Title"
$wsAddUser.Range("B1").Value = "Add user - las"
$wsAddUser.Range("B2").Value = "This field is required."
$wsAddUser.Range("A2").Value = "EmailFieldErrorMessage"

$wsEditUser = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsAddUser)
$wsEditUser.Name = "EditUserPage"
$wsEditUser.Range("A1").Value = "Title"
$wsEditUser.Range("B1").Value = "Edit user - las"

# --- Column widths (Excel ColumnWidth = OOXML width - 0.8333333333333333) -
$wsUsers.Columns.Item(1).ColumnWidth   = 29.830729166666668   # -> 30.6640625
$wsAddUser.Columns.Item(1).ColumnWidth = 27.276041666666668   # -> 28.109375
$wsAddUser.Columns.Item(2).ColumnWidth = 21.830729166666668   # -> 22.6640625
$wsEditUser.Columns.Item(2).ColumnWidth = 16.830729166666668  # -> 17.6640625

# --- Selections on sheets that stay inactive must be set *without*
#     leaving them as the active sheet, so order matters: set the
#     final-active sheet ("EditUserPage") last.
[void]$wsLogin.Range("D25").Select()
[void]$wsHome.Range("B2").Select()
[void]$wsUsers.Range("A4").Select()
[void]$wsAddUser.Range("A1:B1").Select()

# EditUserPage is the sheet left active/selected when the workbook was saved.
[void]$wsEditUser.Activate()
[void]$wsEditUser.Range("O24").Select()
